$d = $word.ActiveDocument

# --- 1. Update date (paragraph 2) ---
$d.Paragraphs.Item(2).Range.Text = 'Updated: 2026-02-18 (local)'

# --- 2. Section heading (paragraph 12) ---
$d.Paragraphs.Item(12).Range.Text = 'New in this update (Railway frontend dependency fix)'

# --- 3-7. Rewrite the five "Railway frontend fix" bullet paragraphs (13-17) ---
$d.Paragraphs.Item(13).Range.Text = '- Fixed frontend Docker build failure due to peer dependency conflict:'
$d.Paragraphs.Item(14).Range.Text = '  - `react-day-picker@8.10.1` expects `date-fns ^2/^3`, while project uses `date-fns@4`.'
$d.Paragraphs.Item(15).Range.Text = '- Updated `frontend/Dockerfile` install command to:'
$d.Paragraphs.Item(16).Range.Text = '  - `npm install --legacy-peer-deps`'
$d.Paragraphs.Item(17).Range.Text = '- This resolves Railway build error `ERESOLVE unable to resolve dependency tree`.'

# --- 8. Insert a new bullet before "- Store email/phone/address ..." (paragraph 20) ---
$storeP = $d.Paragraphs.Item(20)
$storeP.Range.InsertParagraphBefore()
# NOTE: after InsertParagraphBefore(), the original $storeP object reseats itself onto
# the newly-created (empty) paragraph rather than following the "Store email..." text,
# so the new paragraph must be re-fetched fresh by index.
$newP = $d.Paragraphs.Item(20)
$newP.Range.Text = '- Proper long-term dependency alignment still pending (recommended: migrate `date-fns` to a version compatible with all peers or upgrade dependent packages).'

# --- 9. Commit hash (shifted from 26 to 27 after the insertion) ---
$d.Paragraphs.Item(27).Range.Text = '- Last pushed commit: 9cfabad'

# --- 10. Final status line (shifted from 27 to 28) ---
$d.Paragraphs.Item(28).Range.Text = '- Current frontend dependency build fix is local and not pushed yet.'
